$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra (controller-less) field list that used to live in
# column D alongside the "karyawan" block (rows 21-50). This also lets
# Excel prune the now-unreferenced shared strings (idkaryawan, kode_pos,
# kode_pos_domisili, ktp, fotoLama) and shrink the used range/dimension.
$ws.Range("D21:D50").ClearContents()

# Restore the view/selection state recorded for the sheet after the edit.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D20").Select()
